# Add a new "Merge Intervals" entry to the Intervals section of Sheet1,
# and backfill the "Revise 1" date on the existing "Summary Ranges" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 ("Summary Ranges") picks up a revise-on date, matching the date
# formatting already used in column E (numFmtId 17 / "mmm-yy").
$ws.Range("E20").Value = 44621
$ws.Range("E20").NumberFormat = "mmm-yy"

# New row 21: another entry under the "Intervals" category.
$ws.Range("A21").Value = 56
$ws.Range("B21").Value = "Merge Intetvals"
$ws.Range("C21").Value = "Medium"

# Leave the selection where Excel would land after typing the new row.
$ws.Range("C22").Select()
